# Refresh the crypto price/volume snapshot cells (row 2-51) to match
# the latest scrape. Price-looking strings get a leading "'" so Excel
# keeps storing them as text instead of re-parsing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.615.49"
$ws.Range("E2").Value = "  -1.28%  "

$ws.Range("D3").Value = "1.847.79"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'314.61"
$ws.Range("E5").Value = "  -1.13%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4262"
$ws.Range("E7").Value = "  -2.37%  "

$ws.Range("D8").Value = "'0.3651"
$ws.Range("E8").Value = "  -2.31%  "

$ws.Range("D9").Value = "'44.65"
$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("D10").Value = "'0.07314"
$ws.Range("E10").Value = "  -2.29%  "

$ws.Range("D11").Value = "'0.8844"
$ws.Range("E11").Value = "  -5.60%  "

$ws.Range("D12").Value = "'20.77"
$ws.Range("E12").Value = "  -2.51%  "

$ws.Range("D13").Value = "1.872.36"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("D14").Value = "'5.352"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").Value = "'6.545"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("D16").Value = "'0.06923"
$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").Value = "'78.97"
$ws.Range("E18").Value = "  -2.99%  "

$ws.Range("D19").Value = "'0.000008893"
$ws.Range("E19").Value = "  -1.62%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").Value = "'15.45"
$ws.Range("E21").Value = "  -2.82%  "

$ws.Range("D22").Value = "27.624.40"
$ws.Range("E22").Value = "  -1.16%  "

$ws.Range("D23").Value = "'4.984"
$ws.Range("E23").Value = "  -2.76%  "

$ws.Range("D24").Value = "'10.67"
$ws.Range("E24").Value = "  -3.34%  "

$ws.Range("D25").Value = "2.102.78"
$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("D26").Value = "'1.971"
$ws.Range("E26").Value = "  -3.49%  "

$ws.Range("D27").Value = "'153.69"
$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").Value = "'18.98"
$ws.Range("E28").Value = "  +2.57%  "

$ws.Range("D29").Value = "'122.02"
$ws.Range("E29").Value = "  +7.28%  "

$ws.Range("D30").Value = "'5.260"
$ws.Range("E30").Value = "  -5.29%  "

$ws.Range("D31").Value = "'1.917"
$ws.Range("E31").Value = "  +12.31%  "

$ws.Range("D32").Value = "'0.08941"
$ws.Range("E32").Value = "  -0.92%  "

$ws.Range("D33").Value = "'0.7631"
$ws.Range("E33").Value = "  -6.64%  "

$ws.Range("D34").Value = "'4.579"
$ws.Range("E34").Value = "  -5.03%  "

$ws.Range("D35").Value = "'2.979"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").Value = "'1.103"
$ws.Range("E36").Value = "  -6.40%  "

$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D40").Value = "'0.01949"
$ws.Range("E40").Value = "  -1.37%  "

$ws.Range("D41").Value = "'2.818"
$ws.Range("E41").Value = "  -4.89%  "

$ws.Range("D42").Value = "'6.933"
$ws.Range("E42").Value = "  -1.17%  "

$ws.Range("D43").Value = "'0.5119"
$ws.Range("E43").Value = "  -2.88%  "

$ws.Range("D44").Value = "'0.1656"
$ws.Range("E44").Value = "  -2.73%  "

$ws.Range("D45").Value = "'8.282"
$ws.Range("E45").Value = "  -5.79%  "

$ws.Range("D46").Value = "'0.06576"
$ws.Range("E46").Value = "  -2.65%  "

$ws.Range("D47").Value = "'0.4765"
$ws.Range("E47").Value = "  -2.67%  "

$ws.Range("D48").Value = "'10.39"
$ws.Range("E48").Value = "  -2.19%  "

$ws.Range("D49").Value = "'104.39"
$ws.Range("E49").Value = "  -2.55%  "

$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").Value = "'1.631"
$ws.Range("E51").Value = "  -2.64%  "

# Rows 38/39 also swapped coins (Hedera <-> TrustWalletToken), so refresh
# the Coin/Link columns too, not just Price/Volume.
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.096"
$ws.Range("E38").Value = "  -2.42%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05376"
$ws.Range("E39").Value = "  -2.62%  "

